$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 1000
$ws.Cells.Item(76, 9).Value = 1000
$ws.Cells.Item(76, 11).Value = 1000
$ws.Cells.Item(76, 13).Value = -685

$ws.Cells.Item(79, 8).Value = 1000
$ws.Cells.Item(79, 9).Value = 1000
$ws.Cells.Item(79, 11).Value = 1000
$ws.Cells.Item(79, 13).Value = 92

$ws.Cells.Item(98, 8).Value = 8027.9
$ws.Cells.Item(98, 9).Value = 8027.9
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 8027.9
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = -6529.9
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(106, 8).Value = 5076.9287
$ws.Cells.Item(106, 9).Value = 4889.8335
$ws.Cells.Item(106, 11).Value = 4889.8335
$ws.Cells.Item(106, 13).Value = -4258.8335

$ws.Cells.Item(112, 8).Value = 3837.5454
$ws.Cells.Item(112, 10).Value = 3837.5454
$ws.Cells.Item(112, 12).Value = 11512.6362
$ws.Cells.Item(112, 14).Value = -13728.6362

$ws.Cells.Item(122, 8).Value = 8027.9
$ws.Cells.Item(122, 9).Value = 8027.9
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 24083.7
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -21633.7
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(125, 8).Value = 250000270
$ws.Cells.Item(125, 10).Value = 536
$ws.Cells.Item(125, 12).Value = 4824
$ws.Cells.Item(125, 14).Value = -9744

$ws.Cells.Item(129, 8).Value = 2164.4167
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 2164.4167
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 14).Value = -16493.2501
$ws.Cells.Item(129, 12).Value = 6493.250100000001
$ws.Cells.Item(129, 13).ClearContents()

$ws.Cells.Item(135, 8).Value = 2326092.8
$ws.Cells.Item(135, 9).Value = 2381404.5
$ws.Cells.Item(135, 11).Value = 21432640.5
$ws.Cells.Item(135, 13).Value = -21430105.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6254.8887
$ws.Cells.Item(132, 9).Value = 3909.3635
$ws.Cells.Item(132, 11).Value = 11728.0905
$ws.Cells.Item(132, 13).Value = -9198.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7252829.5
$ws.Cells.Item(20, 9).Value = 9263810
$ws.Cells.Item(20, 10).Value = 13298.6
$ws.Cells.Item(20, 11).Value = 9263810
$ws.Cells.Item(20, 12).Value = 13298.6
$ws.Cells.Item(20, 13).Value = -9263563
$ws.Cells.Item(20, 14).Value = -13792.6

$ws.Cells.Item(81, 8).Value = 51659.375
$ws.Cells.Item(81, 10).Value = 51659.375
$ws.Cells.Item(81, 12).Value = 51659.375
$ws.Cells.Item(81, 14).Value = -53781.375

$ws.Cells.Item(84, 8).Value = 51659.375
$ws.Cells.Item(84, 10).Value = 51659.375
$ws.Cells.Item(84, 12).Value = 154978.125
$ws.Cells.Item(84, 14).Value = -165586.125

$ws.Cells.Item(86, 8).Value = 16670150
$ws.Cells.Item(86, 9).Value = 25003080
$ws.Cells.Item(86, 10).Value = 4288.6
$ws.Cells.Item(86, 11).Value = 25003080
$ws.Cells.Item(86, 12).Value = 4288.6
$ws.Cells.Item(86, 13).Value = -25001957
$ws.Cells.Item(86, 14).Value = -6534.6

$ws.Cells.Item(89, 8).Value = 16670150
$ws.Cells.Item(89, 9).Value = 25003080
$ws.Cells.Item(89, 10).Value = 4288.6
$ws.Cells.Item(89, 11).Value = 125015400
$ws.Cells.Item(89, 12).Value = 21443
$ws.Cells.Item(89, 13).Value = -125009784
$ws.Cells.Item(89, 14).Value = -32675

$ws.Cells.Item(105, 8).Value = 3340.6
$ws.Cells.Item(105, 10).Value = 4204.2
$ws.Cells.Item(105, 12).Value = 4204.2
$ws.Cells.Item(105, 14).Value = -7698.2

$ws.Cells.Item(134, 8).Value = 7820374
$ws.Cells.Item(134, 9).Value = 19235088
$ws.Cells.Item(134, 11).Value = 57705264
$ws.Cells.Item(134, 13).Value = -57702729

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8926.296
$ws.Cells.Item(31, 9).Value = 3234.4666
$ws.Cells.Item(31, 10).Value = 16041.083
$ws.Cells.Item(31, 11).Value = 3234.4666
$ws.Cells.Item(31, 12).Value = 16041.083
$ws.Cells.Item(31, 13).Value = -2939.4666
$ws.Cells.Item(31, 14).Value = -16631.083

$ws.Cells.Item(34, 8).Value = 8926.296
$ws.Cells.Item(34, 9).Value = 3234.4666
$ws.Cells.Item(34, 10).Value = 16041.083
$ws.Cells.Item(34, 11).Value = 3234.4666
$ws.Cells.Item(34, 12).Value = 16041.083
$ws.Cells.Item(34, 13).Value = -3032.4666
$ws.Cells.Item(34, 14).Value = -16445.083

$ws.Cells.Item(99, 8).Value = 7021
$ws.Cells.Item(99, 9).Value = 6458.75
$ws.Cells.Item(99, 10).Value = 7583.25
$ws.Cells.Item(99, 11).Value = 6458.75
$ws.Cells.Item(99, 12).Value = 7583.25
$ws.Cells.Item(99, 13).Value = -4960.75
$ws.Cells.Item(99, 14).Value = -10579.25

$ws.Cells.Item(126, 8).Value = 7021
$ws.Cells.Item(126, 9).Value = 6458.75
$ws.Cells.Item(126, 10).Value = 7583.25
$ws.Cells.Item(126, 11).Value = 19376.25
$ws.Cells.Item(126, 12).Value = 22749.75
$ws.Cells.Item(126, 13).Value = -16906.25
$ws.Cells.Item(126, 14).Value = -27689.75

$ws.Cells.Item(132, 8).Value = 4366.655
$ws.Cells.Item(132, 10).Value = 9590.666999999999
$ws.Cells.Item(132, 12).Value = 28772.001
$ws.Cells.Item(132, 14).Value = -33832.001

$ws.Cells.Item(141, 8).Value = 363477.5
$ws.Cells.Item(141, 9).Value = 50000
$ws.Cells.Item(141, 10).Value = 398308.34
$ws.Cells.Item(141, 11).Value = 50000
$ws.Cells.Item(141, 12).Value = 398308.34
$ws.Cells.Item(141, 13).Value = -44820
$ws.Cells.Item(141, 14).Value = -408668.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 576.25
$ws.Cells.Item(103, 9).Value = 320.8
$ws.Cells.Item(103, 10).Value = 1002
$ws.Cells.Item(103, 11).Value = 962.4000000000001
$ws.Cells.Item(103, 12).Value = 3006
$ws.Cells.Item(103, 13).Value = -83.40000000000009
$ws.Cells.Item(103, 14).Value = -4764

$ws.Cells.Item(122, 8).Value = 2016442.1
$ws.Cells.Item(122, 9).Value = 3537388.5
$ws.Cells.Item(122, 10).Value = 910299.25
$ws.Cells.Item(122, 11).Value = 31836496.5
$ws.Cells.Item(122, 12).Value = 8192693.25
$ws.Cells.Item(122, 13).Value = -31834046.5
$ws.Cells.Item(122, 14).Value = -8197593.25

$ws.Cells.Item(129, 8).Value = 111112376
$ws.Cells.Item(129, 9).Value = 2030
$ws.Cells.Item(129, 10).Value = 125001170
$ws.Cells.Item(129, 11).Value = 6090
$ws.Cells.Item(129, 12).Value = 375003510
$ws.Cells.Item(129, 13).Value = -1090
$ws.Cells.Item(129, 14).Value = -375013510

$ws.Cells.Item(131, 8).Value = 2419.9412
$ws.Cells.Item(131, 9).Value = 4343.25
$ws.Cells.Item(131, 11).Value = 13029.75
$ws.Cells.Item(131, 13).Value = -7989.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 71333
$ws.Cells.Item(70, 9).Value = 189999.5
$ws.Cells.Item(70, 11).Value = 189999.5
$ws.Cells.Item(70, 13).Value = -189729.5

$ws.Cells.Item(73, 8).Value = 71333
$ws.Cells.Item(73, 9).Value = 189999.5
$ws.Cells.Item(73, 11).Value = 189999.5
$ws.Cells.Item(73, 13).Value = -189063.5

$ws.Cells.Item(97, 8).Value = 2628.4443
$ws.Cells.Item(97, 9).Value = 2620.8572
$ws.Cells.Item(97, 10).Value = 2655
$ws.Cells.Item(97, 11).Value = 2620.8572
$ws.Cells.Item(97, 12).Value = 2655
$ws.Cells.Item(97, 13).Value = -2124.8572
$ws.Cells.Item(97, 14).Value = -3647

$ws.Cells.Item(122, 8).Value = 10208187
$ws.Cells.Item(122, 9).Value = 35716784
$ws.Cells.Item(122, 10).Value = 4748.8
$ws.Cells.Item(122, 11).Value = 107150352
$ws.Cells.Item(122, 12).Value = 14246.4
$ws.Cells.Item(122, 13).Value = -107147902
$ws.Cells.Item(122, 14).Value = -19146.4

$ws.Cells.Item(126, 8).Value = 7617.8125
$ws.Cells.Item(126, 9).Value = 4822.7334
$ws.Cells.Item(126, 10).Value = 8888.303
$ws.Cells.Item(126, 11).Value = 14468.2002
$ws.Cells.Item(126, 12).Value = 26664.909
$ws.Cells.Item(126, 13).Value = -11998.2002
$ws.Cells.Item(126, 14).Value = -31604.909

$ws.Cells.Item(132, 8).Value = 9796.625
$ws.Cells.Item(132, 9).Value = 6869.75
$ws.Cells.Item(132, 10).Value = 12723.5
$ws.Cells.Item(132, 11).Value = 20609.25
$ws.Cells.Item(132, 12).Value = 38170.5
$ws.Cells.Item(132, 13).Value = -18079.25
$ws.Cells.Item(132, 14).Value = -43230.5

$ws.Cells.Item(133, 8).Value = 90390
$ws.Cells.Item(133, 10).Value = 90390
$ws.Cells.Item(133, 12).Value = 90390
$ws.Cells.Item(133, 14).Value = -100510

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6406.5557
$ws.Cells.Item(7, 9).Value = 5983.143
$ws.Cells.Item(7, 11).Value = 5983.143
$ws.Cells.Item(7, 13).Value = -5871.143

$ws.Cells.Item(22, 8).Value = 3538.8333
$ws.Cells.Item(22, 9).Value = 2666.6667
$ws.Cells.Item(22, 11).Value = 2666.6667
$ws.Cells.Item(22, 13).Value = -2371.6667

$ws.Cells.Item(27, 8).Value = 3538.8333
$ws.Cells.Item(27, 9).Value = 2666.6667
$ws.Cells.Item(27, 11).Value = 2666.6667
$ws.Cells.Item(27, 13).Value = -2559.6667

$ws.Cells.Item(40, 8).Value = 4748.4165
$ws.Cells.Item(40, 9).Value = 3629.5789
$ws.Cells.Item(40, 10).Value = 9000
$ws.Cells.Item(40, 11).Value = 3629.5789
$ws.Cells.Item(40, 12).Value = 9000
$ws.Cells.Item(40, 13).Value = -3493.5789
$ws.Cells.Item(40, 14).Value = -9272

$ws.Cells.Item(122, 8).Value = 6417.2104
$ws.Cells.Item(122, 9).Value = 6126.5713
$ws.Cells.Item(122, 10).Value = 6586.75
$ws.Cells.Item(122, 11).Value = 18379.7139
$ws.Cells.Item(122, 12).Value = 19760.25
$ws.Cells.Item(122, 13).Value = -15929.7139
$ws.Cells.Item(122, 14).Value = -24660.25

$ws.Cells.Item(126, 8).Value = 6406.5557
$ws.Cells.Item(126, 9).Value = 5983.143
$ws.Cells.Item(126, 11).Value = 17949.429
$ws.Cells.Item(126, 13).Value = -15479.429

$ws.Cells.Item(132, 8).Value = 10645814
$ws.Cells.Item(132, 9).Value = 19236060
$ws.Cells.Item(132, 10).Value = 10271.381
$ws.Cells.Item(132, 11).Value = 57708180
$ws.Cells.Item(132, 12).Value = 30814.143
$ws.Cells.Item(132, 13).Value = -57705650
$ws.Cells.Item(132, 14).Value = -35874.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3264.3103
$ws.Cells.Item(122, 9).Value = 2081.875
$ws.Cells.Item(122, 10).Value = 8940
$ws.Cells.Item(122, 11).Value = 6245.625
$ws.Cells.Item(122, 12).Value = 26820
$ws.Cells.Item(122, 13).Value = -3795.625
$ws.Cells.Item(122, 14).Value = -31720

$ws.Cells.Item(126, 8).Value = 1899.9
$ws.Cells.Item(126, 9).Value = 1878.5
$ws.Cells.Item(126, 10).Value = 1949.8334
$ws.Cells.Item(126, 11).Value = 5635.5
$ws.Cells.Item(126, 12).Value = 5849.5002
$ws.Cells.Item(126, 13).Value = -3165.5
$ws.Cells.Item(126, 14).Value = -10789.5002

$ws.Cells.Item(132, 8).Value = 12718.583
$ws.Cells.Item(132, 9).Value = 21709.363
$ws.Cells.Item(132, 10).Value = 5111
$ws.Cells.Item(132, 11).Value = 65128.08900000001
$ws.Cells.Item(132, 12).Value = 15333
$ws.Cells.Item(132, 13).Value = -62598.08900000001
$ws.Cells.Item(132, 14).Value = -20393

$ws.Cells.Item(136, 8).Value = 13703147
$ws.Cells.Item(136, 9).Value = 22225424
$ws.Cells.Item(136, 11).Value = 66676272
$ws.Cells.Item(136, 13).Value = -66673722
